$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New tasks added to the planning (rows 18-23), mirroring the time-format
# used by the existing rows (D/E columns, style copied via NumberFormat "h:mm").
$ws.Range("A18").Value = "Faire l'inscription"
$ws.Range("E18").Value = 0.020833333333333332
$ws.Range("E18").NumberFormat = "h:mm"

$ws.Range("A19").Value = "Faire la dé/connexion"
$ws.Range("E19").Value = 0.013888888888888888
$ws.Range("E19").NumberFormat = "h:mm"

$ws.Range("A20").Value = "Faire la page profil"
$ws.Range("E20").Value = 0.013888888888888888
$ws.Range("E20").NumberFormat = "h:mm"

$ws.Range("A21").Value = "Faire les liens de l'index"
$ws.Range("E21").Value = 0.0069444444444444441
$ws.Range("E21").NumberFormat = "h:mm"

$ws.Range("A22").Value = "Ajouter le bootstrap sur les pages déjà créées"
$ws.Range("D22").NumberFormat = "h:mm"
$ws.Range("E22").Value = 0.076388888888888895
$ws.Range("E22").NumberFormat = "h:mm"

$ws.Range("A23").Value = "Créer le fichier phptohtml"
$ws.Range("E23").Value = 0.020833333333333332
$ws.Range("E23").NumberFormat = "h:mm"

# Update the active selection to match the edited workbook's cursor position.
$ws.Range("E23").Select()
